$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of edits matches the order new shared strings were introduced
# (X,Y,Z -> SPDA -> SPPO -> PING -> Pong's the client. -> Name, Color, Team -> Moves... -> Sets...)
$ws.Range("D5").Value = "X, Y, Z"
$ws.Range("B4").Value = "SPDA"
$ws.Range("B5").Value = "SPPO"
$ws.Range("B3").Value = "PING"
$ws.Range("E3").Value = "Pong's the client."
$ws.Range("D4").Value = "Name, Color, Team"
$ws.Range("E5").Value = "Moves the player to a new position"
$ws.Range("E4").Value = "Sets the player's properties"

# Column widths (best-fit auto resize after the new, wider content was entered)
$ws.Columns.Item(4).ColumnWidth = 17.333333333333336
$ws.Columns.Item(5).ColumnWidth = 31.833333333333336

# Selection
$ws.Range("E5").Select()
